# Auto-generated Word COM-interop script implementing the target diff.
$d = $word.ActiveDocument

function Set-ParagraphRunsXml($paraIndex, $innerRunsXmlB64) {
    $p = $d.Paragraphs($paraIndex)
    $r = $p.Range
    $target = $d.Range($r.Start, $r.End - 1)
    $bytes = [System.Convert]::FromBase64String($innerRunsXmlB64)
    $innerXml = [System.Text.Encoding]::UTF8.GetString($bytes)
    $pkgBytes = [System.Convert]::FromBase64String($script:PkgPrefixB64)
    $pkgPrefix = [System.Text.Encoding]::UTF8.GetString($pkgBytes)
    $sufBytes = [System.Convert]::FromBase64String($script:PkgSuffixB64)
    $pkgSuffix = [System.Text.Encoding]::UTF8.GetString($sufBytes)
    $fullXml = $pkgPrefix + $innerXml + $pkgSuffix
    $target.InsertXML($fullXml)
}

$script:PkgPrefixB64 = "PD94bWwgdmVyc2lvbj0iMS4wIiBlbmNvZGluZz0iVVRGLTgiIHN0YW5kYWxvbmU9InllcyI/Pjxwa2c6cGFja2FnZSB4bWxuczpwa2c9Imh0dHA6Ly9zY2hlbWFzLm1pY3Jvc29mdC5jb20vb2ZmaWNlLzIwMDYveG1sUGFja2FnZSI+PHBrZzpwYXJ0IHBrZzpuYW1lPSIvd29yZC9kb2N1bWVudC54bWwiIHBrZzpjb250ZW50VHlwZT0iYXBwbGljYXRpb24vdm5kLm9wZW54bWxmb3JtYXRzLW9mZmljZWRvY3VtZW50LndvcmRwcm9jZXNzaW5nbWwuZG9jdW1lbnQubWFpbit4bWwiPjxwa2c6eG1sRGF0YT48dzpkb2N1bWVudCB4bWxuczp3PSJodHRwOi8vc2NoZW1hcy5vcGVueG1sZm9ybWF0cy5vcmcvd29yZHByb2Nlc3NpbmdtbC8yMDA2L21haW4iPjx3OmJvZHk+PHc6cD4="
$script:PkgSuffixB64 = "PC93OnA+PC93OmJvZHk+PC93OmRvY3VtZW50PjwvcGtnOnhtbERhdGE+PC9wa2c6cGFydD48L3BrZzpwYWNrYWdlPg=="

Set-ParagraphRunsXml 3 "PHc6cj48dzp0IHhtbDpzcGFjZT0icHJlc2VydmUiPkFqb3V0IGQndW5lIGNvbG9ubmUgZGFucyBsZSB0YWJsZWF1IGRlIHLDqXN1bHRhdHMgcG91ciBpbmRpcXVlciBsZXMgPC93OnQ+PC93OnI+PHc6cj48dzp0PnN0YXR1dHM8L3c6dD48L3c6cj48dzpyPjx3OnQgeG1sOnNwYWNlPSJwcmVzZXJ2ZSI+IChGSU5JU0hFRCwgU1RBUlRFRCwg4oCmKSBhdmVjIGRlcyBpY29uZXMuPC93OnQ+PC93OnI+"
Set-ParagraphRunsXml 12 "PHc6cj48dzp0PkFqb3V0PC93OnQ+PC93OnI+PHc6cj48dzp0IHhtbDpzcGFjZT0icHJlc2VydmUiPiBkdSBuaXZlYXUgZGVzIHByaW9yaXTDqXMgcG91ciBsZXMgPC93OnQ+PC93OnI+PHc6cj48dzp0PnN0YXR1dHM8L3c6dD48L3c6cj48dzpyPjx3OnQgeG1sOnNwYWNlPSJwcmVzZXJ2ZSI+LCA8L3c6dD48L3c6cj48dzpyPjx3OnQ+YXByw6hzPC93OnQ+PC93OnI+PHc6cj48dzp0IHhtbDpzcGFjZT0icHJlc2VydmUiPiBsYSBjb2xsZWN0ZSBkZXMgbWVzc2FnZXMsIHBvdXIgbCdlbnZvaWUgZGFucyBFbGFzdGljU2VhcmNoPC93OnQ+PC93OnI+"
Set-ParagraphRunsXml 13 "PHc6cj48dzp0IHhtbDpzcGFjZT0icHJlc2VydmUiPlRyaSBkZXMgPC93OnQ+PC93OnI+PHc6cj48dzp0PnN0YXR1dHM8L3c6dD48L3c6cj48dzpyPjx3OnQgeG1sOnNwYWNlPSJwcmVzZXJ2ZSI+IHBhciBsZXVyIHByaW9yaXTDqTwvdzp0PjwvdzpyPjx3OnI+PHc6dCB4bWw6c3BhY2U9InByZXNlcnZlIj4gOiBwZXRpdCA8L3c6dD48L3c6cj48dzpyPjx3OnQ+cHJvYmzDqG1lPC93OnQ+PC93OnI+PHc6cj48dzp0PjogbGUgbmc8L3c6dD48L3c6cj48dzpyPjx3OnQ+LTwvdzp0PjwvdzpyPjx3OnI+PHc6dCB4bWw6c3BhY2U9InByZXNlcnZlIj5ncmlkIG5lIHJhZnJhaWNoaSBwYXMgbGVzIGRpcmVjdGl2ZXMgZGUgPC93OnQ+PC93OnI+PHc6cj48dzp0PnN0YXR1dHM8L3c6dD48L3c6cj48dzpyPjx3OnQgeG1sOnNwYWNlPSJwcmVzZXJ2ZSI+LiA8L3c6dD48L3c6cj48dzpyPjx3OnQ+UsOpc29sdXRpb248L3c6dD48L3c6cj48dzpyPjx3OnQgeG1sOnNwYWNlPSJwcmVzZXJ2ZSI+IDogd2F0Y2hlciBzdXIgbGEgY2VsbHVsZSBkZSBsYSBjb2xvbm5lLCBldCBzaSBlbGxlIGNoYW5nZSwgb24gbWV0IGEgam91ciBsJzwvdzp0PjwvdzpyPjx3OnI+PHc6dD5pY29uZTwvdzp0PjwvdzpyPjx3OnI+PHc6dD4sIGxhIGNvdWxldXIgZXQgbGUgdG9vbHRpcC48L3c6dD48L3c6cj4="
Set-ParagraphRunsXml 15 "PHc6cj48dzp0PkTDqWJ1dCBkJ2Fzc29jaWF0aW9uIGQnb2YtY29yZSBldCBvZi08L3c6dD48L3c6cj48dzpyPjx3OnQ+ZnJvbnQuIEplIGRvaXMgbGFuY2VyIGVsYXN0aWNzZWFyY2ggZXQgb2YgY29yZSAoY29sbGVjdGUpIGRlcHVpcyBvZiBmcm9udCwgZGVzIGxlIGxhbmNlbWVudCBkdSBzZXJ2ZXVyIHRvbWNhdC48L3c6dD48L3c6cj4="
Set-ParagraphRunsXml 23 "PHc6cj48dzp0Pklkw6llIGRlIHLDqTwvdzp0PjwvdzpyPjx3OnI+PHc6dD5zbzwvdzp0PjwvdzpyPjx3OnI+PHc6dD5sdXRpb248L3c6dD48L3c6cj48dzpyPjx3OnQgeG1sOnNwYWNlPSJwcmVzZXJ2ZSI+OiA8L3c6dD48L3c6cj48dzpyPjx3OnQ+Y3LDqWF0aW9uPC93OnQ+PC93OnI+PHc6cj48dzp0IHhtbDpzcGFjZT0icHJlc2VydmUiPiBkdSA8L3c6dD48L3c6cj48dzpyPjx3OnQ+Q3JzZkZpbHRlcnM8L3c6dD48L3c6cj48dzpyPjx3OnQgeG1sOnNwYWNlPSJwcmVzZXJ2ZSI+IHBlcnNvIHBvdXIgcXUnaWwgYWNjZXB0ZSBsZXMgPC93OnQ+PC93OnI+PHc6cj48dzp0PnJlcXXDqnRlczwvdzp0PjwvdzpyPjx3OnI+PHc6dCB4bWw6c3BhY2U9InByZXNlcnZlIj4gUE9TVC4gPC93OnQ+PC93OnI+PHc6cj48dzpyUHI+PHc6bGFuZyB3OnZhbD0iZW4tVVMiLz48L3c6clByPjx3OnQgeG1sOnNwYWNlPSJwcmVzZXJ2ZSI+RGVmaW5pdGlvbiBkdSBjc3JmIGZpbHRlciA8L3c6dD48L3c6cj48dzpyPjx3OnQ+ZGFuczwvdzp0PjwvdzpyPjx3OnI+PHc6clByPjx3Omxhbmcgdzp2YWw9ImVuLVVTIi8+PC93OnJQcj48dzp0IHhtbDpzcGFjZT0icHJlc2VydmUiPiBzcHJpbmcgc2VjdXJpdHkgYXZlYyAmbHQ7c2VjdXJpdHk6Y3NyZiByZXF1ZXN0LW1hdGNoZXI8L3c6dD48L3c6cj48dzpyPjx3OnJQcj48dzpsYW5nIHc6dmFsPSJlbi1VUyIvPjwvdzpyUHI+PHc6dD4tcmVmPC93OnQ+PC93OnI+PHc6cj48dzpyUHI+PHc6bGFuZyB3OnZhbD0iZW4tVVMiLz48L3c6clByPjx3OnQ+PSJtYWNsYXNzIi8mZ3Q7PC93OnQ+PC93OnI+"
Set-ParagraphRunsXml 24 "PHc6cj48dzp0IHhtbDpzcGFjZT0icHJlc2VydmUiPkZvbmN0aW9ubmUgLSZndDsgYWpvdXQgZGVzIGRvbm7DqWVzIGTDqWrDoCA8L3c6dD48L3c6cj48dzpyPjx3OnQ+Y3LDqWVyPC93OnQ+PC93OnI+PHc6cj48dzp0IHhtbDpzcGFjZT0icHJlc2VydmUiPiBzdXIgbGUgc2VydmV1ciwgcG91ciBsZXMgPC93OnQ+PC93OnI+PHc6cj48dzp0PmludMOpZ3Jlcjwvdzp0PjwvdzpyPjx3OnI+PHc6dCB4bWw6c3BhY2U9InByZXNlcnZlIj4gZGFucyBsZSBkb3NzaWVyIHNww6ljaWZpw6k8L3c6dD48L3c6cj4="
Set-ParagraphRunsXml 25 "PHc6cj48dzp0PlLDqXNvbHV0aW9uPC93OnQ+PC93OnI+PHc6cj48dzp0IHhtbDpzcGFjZT0icHJlc2VydmUiPiBkZXMgZXJyZXVycyBhdSA8L3c6dD48L3c6cj48dzpyPjx3OnQ+ZMOpbWFycmFnZTwvdzp0PjwvdzpyPjx3OnI+PHc6dCB4bWw6c3BhY2U9InByZXNlcnZlIj4gZGUgdG9tY2F0LiBFeCA6IGxvZzRqIHF1aSBuJ2EgcGFzIDwvdzp0PjwvdzpyPjx3OnI+PHc6dD5hY2PDqHM8L3c6dD48L3c6cj48dzpyPjx3OnQgeG1sOnNwYWNlPSJwcmVzZXJ2ZSI+IGF1IGZpY2hpZXIgZGUgbG9nPC93OnQ+PC93OnI+"
Set-ParagraphRunsXml 26 "PHc6cj48dzpsYXN0UmVuZGVyZWRQYWdlQnJlYWsvPjx3OnQ+QW3DqWxpb3Jlcjwvdzp0PjwvdzpyPjx3OnI+PHc6dCB4bWw6c3BhY2U9InByZXNlcnZlIj4gbCdvdXRpbCBkZSB2ZXJzaW9ubmluZyBkZSBkb2N1bWVudHMuPC93OnQ+PC93OnI+PHc6cj48dzp0IHhtbDpzcGFjZT0icHJlc2VydmUiPiAyIGxvZ2lxdWVzIGRlIHZlcnNpb25uaW5nIHByb3Bvc8OpLCA8L3c6dD48L3c6cj48dzpyPjx3OnQ+w6A8L3c6dD48L3c6cj48dzpyPjx3OnQgeG1sOnNwYWNlPSJwcmVzZXJ2ZSI+IG1vaSBkZSBjaG9pc2lyIGwndW5lIGQnZWxsZS4gPC93OnQ+PC93OnI+"
Set-ParagraphRunsXml 27 "PHc6cj48dzp0IHhtbDpzcGFjZT0icHJlc2VydmUiPkxvZ2lxdWUgMSA6IENvbXBhcmVyIGxlIG5vdXZlYXUgZG9jdW1lbnQgKGF2YW50IGluc2VydGlvbiksIGF2ZWMgbGUgZG9jdW1lbnQgYWN0aWYuIFB1aXMgYWpvdXRlciBkYW5zIGxhIGJhc2UgZGUgdmVyc2lvbm5pbmcsIHF1ZSBsZXMgY2hhbXBzIHF1aSBkaWZmw6hyZSBlbnRyZSBsZXMgMi4gRXQgPC93OnQ+PC93OnI+PHc6cj48dzp0Pmluc8OpcmVyPC93OnQ+PC93OnI+PHc6cj48dzp0IHhtbDpzcGFjZT0icHJlc2VydmUiPiBsZSBub3V2ZWF1IGRvY3VtZW50IGNvbW1lIGFjdGlmLjwvdzp0PjwvdzpyPg=="
Set-ParagraphRunsXml 28 "PHc6cj48dzp0IHhtbDpzcGFjZT0icHJlc2VydmUiPkxvZ2lxdWUgMiA6IDwvdzp0PjwvdzpyPjx3OnI+PHc6dD5JbnPDqXJlcjwvdzp0PjwvdzpyPjx3OnI+PHc6dCB4bWw6c3BhY2U9InByZXNlcnZlIj4gZGFucyBsYSBiYXNlIGRlIHZlcnNpb25uaW5nIGxlIGRvY3VtZW50IGFjdGlmIGF1IGNvbXBsZXQsIHB1aXMgbWV0dHJlIGxlIG5vdXZlYXUgZG9jdW1lbnQgY29tbWUgYWN0aWYuIEV0IGZhaXJlIHRvdXJuw6llIHVuZSByb3V0aW5lIHF1aSB2YSBjb21wYXJlciBsZXMgY2hhbXBzIGV0IHN1cHByaW1lciBjZXV4IHF1aSBzb250IGlkZW50aXF1ZXMuPC93OnQ+PC93OnI+"
Set-ParagraphRunsXml 29 "PHc6cj48dzp0PlLDqWZsZXhpb248L3c6dD48L3c6cj48dzpyPjx3OnQgeG1sOnNwYWNlPSJwcmVzZXJ2ZSI+IC0mZ3Q7IGxvZ2lxdWUgMiwgY29tcGxpcXXDqSBldCA8L3c6dD48L3c6cj48dzpyPjx3OnQ+c3VyY2hhcmdlPC93OnQ+PC93OnI+PHc6cj48dzp0IHhtbDpzcGFjZT0icHJlc2VydmUiPiBkZSBib3VjbGUgcG91ciBjb21wYXLDqSBkYW5zIGxhIGJhc2UgZGUgdmVyc2lvbm5pbmcgY2UgcXVpIGEgY2hhbmdlciBldGMuPC93OnQ+PC93OnI+"
Set-ParagraphRunsXml 30 "PHc6cj48dzp0IHhtbDpzcGFjZT0icHJlc2VydmUiPkonYWkgY2hvaXNpIGxhIGxvZ2lxdWUgMSwgPC93OnQ+PC93OnI+PHc6cj48dzp0PmFwcsOoczwvdzp0PjwvdzpyPjx3OnI+PHc6dCB4bWw6c3BhY2U9InByZXNlcnZlIj4gYXZvaXIgPC93OnQ+PC93OnI+PHc6cj48dzp0PnRyb3V2w6k8L3c6dD48L3c6cj48dzpyPjx3OnQgeG1sOnNwYWNlPSJwcmVzZXJ2ZSI+IHVuZSA8L3c6dD48L3c6cj48dzpyPjx3OnQ+QVBJPC93OnQ+PC93OnI+PHc6cj48dzp0IHhtbDpzcGFjZT0icHJlc2VydmUiPiBqYXZhIHF1aSBjb21wYXJlIDIgb2JqZXRzLCBldCByZW52b2llIGxlcyA8L3c6dD48L3c6cj48dzpyPjx3OnQ+ZGlmZsOpcmVuY2VzPC93OnQ+PC93OnI+PHc6cj48dzp0IHhtbDpzcGFjZT0icHJlc2VydmUiPiBlbnRyZSBsZXMgMiwgY2UgcXVpIHBlcm1ldCBkZSBzYXZvaXIgY2UgcXVpIGEgPC93OnQ+PC93OnI+PHc6cj48dzp0PnLDqWVsbGVtZW50PC93OnQ+PC93OnI+PHc6cj48dzp0IHhtbDpzcGFjZT0icHJlc2VydmUiPiBjaGFuZ2VyLCBldCBub3VzIHBlcm1ldHRyYSBkJ2Fqb3V0ZXIgZGlyZWN0ZW1lbnQgZGFucyBsYSBiYXNlIGRlIHZlcnNpb25uaW5nIHF1ZSBsZXMgY2hhbXBzIHF1aSBvbnQgw6l0w6kgbW9kaWZpw6kuPC93OnQ+PC93OnI+"

# Append new paragraph after paragraph 32 ("Implémentation de Javers dans le projet")
$pLast = $d.Paragraphs(32)
$rLast = $pLast.Range
$insertPoint = $d.Range($rLast.End, $rLast.End)
$newParaBytes = [System.Convert]::FromBase64String("PD94bWwgdmVyc2lvbj0iMS4wIiBlbmNvZGluZz0iVVRGLTgiIHN0YW5kYWxvbmU9InllcyI/Pjxwa2c6cGFja2FnZSB4bWxuczpwa2c9Imh0dHA6Ly9zY2hlbWFzLm1pY3Jvc29mdC5jb20vb2ZmaWNlLzIwMDYveG1sUGFja2FnZSI+PHBrZzpwYXJ0IHBrZzpuYW1lPSIvd29yZC9kb2N1bWVudC54bWwiIHBrZzpjb250ZW50VHlwZT0iYXBwbGljYXRpb24vdm5kLm9wZW54bWxmb3JtYXRzLW9mZmljZWRvY3VtZW50LndvcmRwcm9jZXNzaW5nbWwuZG9jdW1lbnQubWFpbit4bWwiPjxwa2c6eG1sRGF0YT48dzpkb2N1bWVudCB4bWxuczp3PSJodHRwOi8vc2NoZW1hcy5vcGVueG1sZm9ybWF0cy5vcmcvd29yZHByb2Nlc3NpbmdtbC8yMDA2L21haW4iPjx3OmJvZHk+PHc6cD48dzpwUHI+PHc6cFN0eWxlIHc6dmFsPSJQYXJhZ3JhcGhlZGVsaXN0ZSIvPjx3Om51bVByPjx3Omlsdmwgdzp2YWw9IjIiLz48dzpudW1JZCB3OnZhbD0iMSIvPjwvdzpudW1Qcj48L3c6cFByPjx3OnI+PHc6dCB4bWw6c3BhY2U9InByZXNlcnZlIj5QYXMgPC93OnQ+PC93OnI+PHc6cHJvb2ZFcnIgdzp0eXBlPSJzcGVsbFN0YXJ0Ii8+PHc6cj48dzp0PmJjcDwvdzp0PjwvdzpyPjx3OnByb29mRXJyIHc6dHlwZT0ic3BlbGxFbmQiLz48dzpyPjx3OnQgeG1sOnNwYWNlPSJwcmVzZXJ2ZSI+IGRlIGRvYywgZGlmZmljaWxlIGQnPC93OnQ+PC93OnI+PHc6cHJvb2ZFcnIgdzp0eXBlPSJzcGVsbFN0YXJ0Ii8+PHc6cj48dzp0PmFjY2VkZXI8L3c6dD48L3c6cj48dzpwcm9vZkVyciB3OnR5cGU9InNwZWxsRW5kIi8+PHc6cj48dzp0IHhtbDpzcGFjZT0icHJlc2VydmUiPiA8L3c6dD48L3c6cj48dzpwcm9vZkVyciB3OnR5cGU9ImdyYW1TdGFydCIvPjx3OnI+PHc6dD5hdTwvdzp0PjwvdzpyPjx3OnByb29mRXJyIHc6dHlwZT0iZ3JhbUVuZCIvPjx3OnI+PHc6dCB4bWw6c3BhY2U9InByZXNlcnZlIj4gcHJvcHJpw6l0w6kgZXQgYSBsZXVyIHZhbHVlIChjZWxsZSBxdWkgb250IMOpdMOpIG1vZGlmacOpKTwvdzp0PjwvdzpyPjwvdzpwPjwvdzpib2R5Pjwvdzpkb2N1bWVudD48L3BrZzp4bWxEYXRhPjwvcGtnOnBhcnQ+PC9wa2c6cGFja2FnZT4=")
$newParaXml = [System.Text.Encoding]::UTF8.GetString($newParaBytes)
$insertPoint.InsertXML($newParaXml)

Write-Host "Done."
